$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells stay as text (matches source inlineStr values),
# since plain numeric-looking strings would otherwise be auto-converted to numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.910.65"
$ws.Range("E2").Value = "  -0.37%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.856.35"
$ws.Range("E3").Value = "  -1.54%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.03"
$ws.Range("E5").Value = "  -0.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.20%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5136"
$ws.Range("E7").Value = "  +2.70%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3803"
$ws.Range("E8").Value = "  -0.53%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08269"
$ws.Range("E9").Value = "  -10.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.59"
$ws.Range("E10").Value = "  -0.19%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.104"
$ws.Range("E11").Value = "  -1.60%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.161"
$ws.Range("E12").Value = "  -2.65%  "
$ws.Range("B13").Value = "Solana"
$ws.Range("C13").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.38"
$ws.Range("E13").Value = "  -1.81%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.853.92"
$ws.Range("E14").Value = "  -1.76%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.176"
$ws.Range("E15").Value = "  -1.37%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.003"
$ws.Range("E16").Value = "  +0.24%  "
$ws.Range("E17").Value = "  -1.17%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "90.19"
$ws.Range("E18").Value = "  -1.32%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06591"
$ws.Range("E19").Value = "  -0.86%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.66"
$ws.Range("E21").Value = "  +0.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.984"
$ws.Range("E22").Value = "  -3.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.958.50"
$ws.Range("E23").Value = "  -0.34%  "
$ws.Range("E24").Value = "  -3.74%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.212"
$ws.Range("E25").Value = "  -3.83%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.564"
$ws.Range("E26").Value = "  +0.55%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.074.42"
$ws.Range("E27").Value = "  -1.44%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "156.61"
$ws.Range("E28").Value = "  -0.51%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.32"
$ws.Range("E29").Value = "  -2.33%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "124.04"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1061"
$ws.Range("E31").Value = "  +0.68%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.034"
$ws.Range("E32").Value = "  -3.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.578"
$ws.Range("E33").Value = "  -0.31%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.591"
$ws.Range("E34").Value = "  +0.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.475"
$ws.Range("E35").Value = "  +1.40%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06509"
$ws.Range("E36").Value = "  -1.01%  "
$ws.Range("E37").Value = "  -0.46%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2167"
$ws.Range("E38").Value = "  -1.46%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.203"
$ws.Range("E39").Value = "  -0.64%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6424"
$ws.Range("E40").Value = "  +0.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.227"
$ws.Range("E41").Value = "  -4.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.21"
$ws.Range("E42").Value = "  -2.96%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.865"
$ws.Range("E43").Value = "  -1.54%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6087"
$ws.Range("E44").Value = "  +0.72%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.97"
$ws.Range("E45").Value = "  -2.67%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.280"
$ws.Range("E46").Value = "  -0.79%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.652"
$ws.Range("E47").Value = "  -0.53%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.966"
$ws.Range("E48").Value = "  -1.28%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.200"
$ws.Range("E49").Value = "  -1.15%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "120.38"
$ws.Range("E50").Value = "  -0.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.37"
$ws.Range("E51").Value = "  +0.60%  "
